$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "is_locked" (O1) and "is_enabled" (P1) columns, shifting the
# trailing "rem" column (previously Q1) left into O1.
$ws.Range("O1:P1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
